$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that vary per data row (2..13) and must be reshuffled.
$cols = @("D","J","K","L","M","N","O","P","Q")

# For each destination row (key), the source row whose data should be copied there.
# (Derived from matching the full row tuples between the old and new sheet state.)
$mapping = @{
    2  = 7
    3  = 13
    4  = 2
    5  = 6
    6  = 3
    7  = 11
    8  = 4
    9  = 10
    10 = 8
    11 = 12
    12 = 9
    13 = 5
}

# Snapshot the current values of the relevant columns for every data row first,
# since rows will be overwritten in place and some are sources for others.
$snapshot = @{}
for ($r = 2; $r -le 13; $r++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Apply the permutation: row $destRow gets the values that used to live in row $mapping[$destRow].
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $srcVals[$col]
    }
}
